# Apply the weekly-report update for WR_89775268_WeekEnding_062925.xlsx
# - refresh the "Report Generated On" timestamp
# - zero out all pricing (now $0), bump the line-item count
# - insert a new "Point 07 / TIE-4-ALH-F" line item on Saturday's table
#   (pushes the remaining Point 09 rows + TOTAL row down by one)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates -------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 12

# --- Thursday table: zero out pricing -----------------------------------------
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0

# --- Saturday table: zero out pricing for existing Point 07 rows --------------
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0

# --- Insert the new "Point 07 / TIE-4-ALH-F" line item at row 28 -------------
# This shifts the old rows 28-32 (Point 09 items + TOTAL row) down to 29-33,
# and shifts the merged TOTAL label range (A32:G32 -> A33:G33) automatically.
$ws.Rows(28).EntireRow.Insert()

# The freshly inserted row inherits formatting from the row above (the odd
# striping); copy formats from the row below instead so the even/odd banding
# stays correct (row 28 should look like the row that used to be row 28, now
# shifted to row 29).
$ws.Range("A29:H29").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A28").Value = "Point 07"
$ws.Range("B28").Value = "TIE-4-ALH-F"
$ws.Range("C28").Value = "Inst"
$ws.Range("D28").Value = "TIE,4 AWG,AL Hand Tie,F Neck"
$ws.Range("E28").Value = "EA"
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = ""
$ws.Range("H28").Value = 0

# --- Zero out pricing on the rows that shifted down (old rows 28-31) ---------
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0

# --- TOTAL row (shifted from 32 to 33) ----------------------------------------
$ws.Range("H33").Value = 0
